$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.022103
$ws.Range("H2").Value = 0.066309
$ws.Range("I2").Value = 0.0007043476645371027
$ws.Range("J2").Value = 0.0007043476645371028
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 3.825035
$ws.Range("N2").Value = 11.475105
$ws.Range("O2").Value = 0.03111562857396839
$ws.Range("P2").Value = 0.03111562857396839
$ws.Range("Q2").Value = 0.084544748605
$ws.Range("R2").Value = 0.760902737445
$ws.Range("S2").Value = 0.00002191622031667857
$ws.Range("T2").Value = 0.00002191622031667858

# Row 3
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.022103
$ws.Range("H3").Value = 0.066309
$ws.Range("I3").Value = 0.0007043476645371027
$ws.Range("J3").Value = 0.0007043476645371028
$ws.Range("O3").Value = 0.4709815605157605
$ws.Range("P3").Value = 0.4709815605157605
$ws.Range("Q3").Value = 1.279711175904333
$ws.Range("R3").Value = 11.517400583139
$ws.Range("S3").Value = 0.000331734762189316
$ws.Range("T3").Value = 0.000331734762189316

# Row 4
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.022103
$ws.Range("H4").Value = 0.066309
$ws.Range("I4").Value = 0.0007043476645371027
$ws.Range("J4").Value = 0.0007043476645371028
$ws.Range("M4").Value = 61.10114166666667
$ws.Range("N4").Value = 183.303425
$ws.Range("O4").Value = 0.4970413158429724
$ws.Range("P4").Value = 0.4970413158429724
$ws.Range("Q4").Value = 1.350518534258333
$ws.Range("R4").Value = 12.154666808325
$ws.Range("S4").Value = 0.000350089889992446
$ws.Range("T4").Value = 0.000350089889992446

# Row 5
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.022103
$ws.Range("H5").Value = 0.066309
$ws.Range("I5").Value = 0.0007043476645371027
$ws.Range("J5").Value = 0.0007043476645371028
$ws.Range("M5").Value = 0.1059033333333333
$ws.Range("N5").Value = 0.31771
$ws.Range("O5").Value = 0.0008614950672987739
$ws.Range("P5").Value = 0.0008614950672987739
$ws.Range("Q5").Value = 0.002340781376666667
$ws.Range("R5").Value = 0.02106703239
$ws.Range("S5").Value = 0.0000006067920386621254
$ws.Range("T5").Value = 0.0000006067920386621255

# Row 6
$ws.Range("G6").Value = 5.827140333333332
$ws.Range("I6").Value = 0.1856912041222136
$ws.Range("J6").Value = 0.1856912041222136
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 3.825035
$ws.Range("N6").Value = 11.475105
$ws.Range("O6").Value = 0.03111562857396839
$ws.Range("P6").Value = 0.03111562857396839
$ws.Range("Q6").Value = 22.28901572491166
$ws.Range("R6").Value = 200.601141524205
$ws.Range("S6").Value = 0.005777898536919745
$ws.Range("T6").Value = 0.005777898536919746

# Row 7
$ws.Range("G7").Value = 5.827140333333332
$ws.Range("I7").Value = 0.1856912041222136
$ws.Range("J7").Value = 0.1856912041222136
$ws.Range("O7").Value = 0.4709815605157605
$ws.Range("P7").Value = 0.4709815605157605
$ws.Range("S7").Value = 0.08745713309153076
$ws.Range("T7").Value = 0.08745713309153076

# Row 8
$ws.Range("G8").Value = 5.827140333333332
$ws.Range("I8").Value = 0.1856912041222136
$ws.Range("J8").Value = 0.1856912041222136
$ws.Range("M8").Value = 61.10114166666667
$ws.Range("N8").Value = 183.303425
$ws.Range("O8").Value = 0.4970413158429724
$ws.Range("P8").Value = 0.4970413158429724
$ws.Range("Q8").Value = 356.0449270185472
$ws.Range("R8").Value = 3204.404343166925
$ws.Range("S8").Value = 0.09229620043737101
$ws.Range("T8").Value = 0.09229620043737101

# Row 9
$ws.Range("G9").Value = 5.827140333333332
$ws.Range("I9").Value = 0.1856912041222136
$ws.Range("J9").Value = 0.1856912041222136
$ws.Range("M9").Value = 0.1059033333333333
$ws.Range("N9").Value = 0.31771
$ws.Range("O9").Value = 0.0008614950672987739
$ws.Range("P9").Value = 0.0008614950672987739
$ws.Range("Q9").Value = 0.617113585101111
$ws.Range("R9").Value = 5.554022265909999
$ws.Range("S9").Value = 0.0001599720563920567
$ws.Range("T9").Value = 0.0001599720563920567

# Row 10
$ws.Range("G10").Value = 15.496839
$ws.Range("H10").Value = 46.490517
$ws.Range("I10").Value = 0.4938317132225258
$ws.Range("J10").Value = 0.4938317132225258
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 3.825035
$ws.Range("N10").Value = 11.475105
$ws.Range("O10").Value = 0.03111562857396839
$ws.Range("P10").Value = 0.03111562857396839
$ws.Range("Q10").Value = 59.27595156436499
$ws.Range("R10").Value = 533.483564079285
$ws.Range("S10").Value = 0.01536588416667859
$ws.Range("T10").Value = 0.01536588416667859

# Row 11
$ws.Range("G11").Value = 15.496839
$ws.Range("H11").Value = 46.490517
$ws.Range("I11").Value = 0.4938317132225258
$ws.Range("J11").Value = 0.4938317132225258
$ws.Range("O11").Value = 0.4709815605157605
$ws.Range("P11").Value = 0.4709815605157605
$ws.Range("Q11").Value = 897.230152444923
$ws.Range("R11").Value = 8075.071372004306
$ws.Range("S11").Value = 0.2325856309257167
$ws.Range("T11").Value = 0.2325856309257167

# Row 12
$ws.Range("G12").Value = 15.496839
$ws.Range("H12").Value = 46.490517
$ws.Range("I12").Value = 0.4938317132225258
$ws.Range("J12").Value = 0.4938317132225258
$ws.Range("M12").Value = 61.10114166666667
$ws.Range("N12").Value = 183.303425
$ws.Range("O12").Value = 0.4970413158429724
$ws.Range("P12").Value = 0.4970413158429724
$ws.Range("Q12").Value = 946.874555124525
$ws.Range("R12").Value = 8521.870996120724
$ws.Range("S12").Value = 0.2454547645451136
$ws.Range("T12").Value = 0.2454547645451136

# Row 13
$ws.Range("G13").Value = 15.496839
$ws.Range("H13").Value = 46.490517
$ws.Range("I13").Value = 0.4938317132225258
$ws.Range("J13").Value = 0.4938317132225258
$ws.Range("M13").Value = 0.1059033333333333
$ws.Range("N13").Value = 0.31771
$ws.Range("O13").Value = 0.0008614950672987739
$ws.Range("P13").Value = 0.0008614950672987739
$ws.Range("Q13").Value = 1.64116690623
$ws.Range("R13").Value = 14.77050215607
$ws.Range("S13").Value = 0.0004254335850169087
$ws.Range("T13").Value = 0.0004254335850169087

# Row 14
$ws.Range("G14").Value = 10.03472733333333
$ws.Range("H14").Value = 30.104182
$ws.Range("I14").Value = 0.3197727349907235
$ws.Range("J14").Value = 0.3197727349907235
$ws.Range("K14").Value = 3
$ws.Range("L14").Value = 1
$ws.Range("M14").Value = 3.825035
$ws.Range("N14").Value = 11.475105
$ws.Range("O14").Value = 0.03111562857396839
$ws.Range("P14").Value = 0.03111562857396839
$ws.Range("Q14").Value = 38.38318326545667
$ws.Range("R14").Value = 345.44864938911
$ws.Range("S14").Value = 0.009949929650053377
$ws.Range("T14").Value = 0.009949929650053377

# Row 15
$ws.Range("G15").Value = 10.03472733333333
$ws.Range("H15").Value = 30.104182
$ws.Range("I15").Value = 0.3197727349907235
$ws.Range("J15").Value = 0.3197727349907235
$ws.Range("O15").Value = 0.4709815605157605
$ws.Range("P15").Value = 0.4709815605157605
$ws.Range("Q15").Value = 580.986866742947
$ws.Range("R15").Value = 5228.881800686522
$ws.Range("S15").Value = 0.1506070617363237
$ws.Range("T15").Value = 0.1506070617363237

# Row 16
$ws.Range("G16").Value = 10.03472733333333
$ws.Range("H16").Value = 30.104182
$ws.Range("I16").Value = 0.3197727349907235
$ws.Range("J16").Value = 0.3197727349907235
$ws.Range("M16").Value = 61.10114166666667
$ws.Range("N16").Value = 183.303425
$ws.Range("O16").Value = 0.4970413158429724
$ws.Range("P16").Value = 0.4970413158429724
$ws.Range("Q16").Value = 613.1332963803724
$ws.Range("R16").Value = 5518.199667423351
$ws.Range("S16").Value = 0.1589402609704953
$ws.Range("T16").Value = 0.1589402609704953

# Row 17
$ws.Range("G17").Value = 10.03472733333333
$ws.Range("H17").Value = 30.104182
$ws.Range("I17").Value = 0.3197727349907235
$ws.Range("J17").Value = 0.3197727349907235
$ws.Range("M17").Value = 0.1059033333333333
$ws.Range("N17").Value = 0.31771
$ws.Range("O17").Value = 0.0008614950672987739
$ws.Range("P17").Value = 0.0008614950672987739
$ws.Range("Q17").Value = 1.062711073691111
$ws.Range("R17").Value = 9.56439966322
$ws.Range("S17").Value = 0.0002754826338511463
$ws.Range("T17").Value = 0.0002754826338511463

